$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 139, shifting existing rows 139-169 down to 145-175
$ws.Rows("139:144").Insert()

# New data for the 6 inserted rows (139-144)
$newRows = @(
    @{ A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=44588; E=4; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103001; J="Cereza"; K="Lapins";  L="Especial"; M=500; N=11500; O=12000; P=11750; Q="$/bandeja 10 kilos"; R="Provincia de Curicó"; S=1175; T=10 },
    @{ A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=44588; E=4; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103001; J="Cereza"; K="Lapins";  L="Primera";  M=400; N=9500;  O=10000; P=9750;  Q="$/bandeja 10 kilos"; R="Provincia de Curicó"; S=975;  T=10 },
    @{ A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=44588; E=4; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103001; J="Cereza"; K="Lapins";  L="Segunda";  M=360; N=7500;  O=8000;  P=7750;  Q="$/bandeja 10 kilos"; R="Provincia de Curicó"; S=775;  T=10 },
    @{ A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=44588; E=4; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103001; J="Cereza"; K="Santina"; L="Especial"; M=400; N=12500; O=13000; P=12750; Q="$/bandeja 10 kilos"; R="Provincia de Curicó"; S=1275; T=10 },
    @{ A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=44588; E=4; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103001; J="Cereza"; K="Santina"; L="Primera";  M=400; N=10500; O=11000; P=10750; Q="$/bandeja 10 kilos"; R="Provincia de Curicó"; S=1075; T=10 },
    @{ A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=44588; E=4; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103001; J="Cereza"; K="Santina"; L="Segunda";  M=360; N=8500;  O=9000;  P=8750;  Q="$/bandeja 10 kilos"; R="Provincia de Curicó"; S=875;  T=10 }
)

$startRow = 139
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = $data.T
}
